# Applies the diff to the workbook:
#  - About sheet: remove the old "Transmission Capacity That Provides Flexibility "
#    shared string, replace B8's inline note text with a multi-paragraph
#    explanation placed in A9:A13, and move the original note text down to A15.
#  - FoTCAMRBtPF sheet: update the header text in B1 to the new wording.

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("FoTCAMRBtPF")

# --- About sheet -----------------------------------------------------

# Row 8 used to hold "Notes:" in A8 and the note text itself in B8.
# The note text moves down to A15, so clear B8 entirely.
$wsAbout.Range("B8").ClearContents()

# The old A9 was an empty, bold-styled placeholder cell - clear both its
# contents and formatting so it becomes a plain, unstyled cell like the
# new rows below it.
$wsAbout.Range("A9").ClearContents()
$wsAbout.Range("A9").ClearFormats()

# Row 14 stays empty; the original note text is relocated to row 15.
$wsAbout.Range("A15").Value = "All imports into India are from hydro facilities, so we assume these are fully flexible."

# --- FoTCAMRBtPF sheet ------------------------------------------------

# Header text updated to include units and drop the trailing space.
$wsData.Range("B1").Value = "Transmission Capacity That Provides Flexibility (dimensionless)"

$wsData.Range("C7").Select()

# --- About sheet (continued) ------------------------------------------

# New explanatory paragraph, one sentence/line per row, unstyled.
$wsAbout.Range("A9").Value  = "This variable represents the share of the transmission capacity across the modeled"
$wsAbout.Range("A10").Value = "region border that can be used to provide flexibility.  This may be due to flexible"
$wsAbout.Range("A11").Value = "generation outside the modeled region, or it may be due to differences in"
$wsAbout.Range("A12").Value = "temporal alignment between inflexible resources within the modeled region and"
$wsAbout.Range("A13").Value = "inflexible resources outside the modeled region."

# Keep the "About" sheet as the active/selected tab, matching the target file,
# and update the selection to match the new layout.
$wsAbout.Activate()
$wsAbout.Range("A9:A15").Select()
